$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("33÷9=3, 6", $true, $false, $false, $false, $false, $true, 0, $false, "46÷6=7, 4", 1) | Out-Null

$cell = $t.Cell(1, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("58÷7=8, 2", $true, $false, $false, $false, $false, $true, 0, $false, "11÷8=1, 3", 1) | Out-Null

$cell = $t.Cell(1, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("31÷2=15, 1", $true, $false, $false, $false, $false, $true, 0, $false, "72÷2=36, 0", 1) | Out-Null

$cell = $t.Cell(1, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("56÷4=14, 0", $true, $false, $false, $false, $false, $true, 0, $false, "90÷6=15, 0", 1) | Out-Null

$cell = $t.Cell(1, 5)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("51÷3=17, 0", $true, $false, $false, $false, $false, $true, 0, $false, "14÷6=2, 2", 1) | Out-Null

$cell = $t.Cell(5, 1)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("58÷4=14, 2", $true, $false, $false, $false, $false, $true, 0, $false, "12÷2=6, 0", 1) | Out-Null

$cell = $t.Cell(5, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("33÷9=3, 6", $true, $false, $false, $false, $false, $true, 0, $false, "88÷6=14, 4", 1) | Out-Null

$cell = $t.Cell(5, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("31÷9=3, 4", $true, $false, $false, $false, $false, $true, 0, $false, "39÷2=19, 1", 1) | Out-Null

$cell = $t.Cell(5, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("64÷5=12, 4", $true, $false, $false, $false, $false, $true, 0, $false, "91÷4=22, 3", 1) | Out-Null

$cell = $t.Cell(5, 5)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("32÷6=5, 2", $true, $false, $false, $false, $false, $true, 0, $false, "94÷6=15, 4", 1) | Out-Null

$cell = $t.Cell(9, 1)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("81÷6=13, 3", $true, $false, $false, $false, $false, $true, 0, $false, "27÷4=6, 3", 1) | Out-Null

$cell = $t.Cell(9, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("28÷3=9, 1", $true, $false, $false, $false, $false, $true, 0, $false, "60÷8=7, 4", 1) | Out-Null

$cell = $t.Cell(9, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("84÷4=21, 0", $true, $false, $false, $false, $false, $true, 0, $false, "43÷4=10, 3", 1) | Out-Null

$cell = $t.Cell(9, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("91÷4=22, 3", $true, $false, $false, $false, $false, $true, 0, $false, "83÷7=11, 6", 1) | Out-Null

$cell = $t.Cell(9, 5)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("27÷5=5, 2", $true, $false, $false, $false, $false, $true, 0, $false, "61÷9=6, 7", 1) | Out-Null

$cell = $t.Cell(13, 1)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("67÷4=16, 3", $true, $false, $false, $false, $false, $true, 0, $false, "61÷7=8, 5", 1) | Out-Null

$cell = $t.Cell(13, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("31÷4=7, 3", $true, $false, $false, $false, $false, $true, 0, $false, "27÷2=13, 1", 1) | Out-Null

$cell = $t.Cell(13, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("44÷2=22, 0", $true, $false, $false, $false, $false, $true, 0, $false, "20÷6=3, 2", 1) | Out-Null

$cell = $t.Cell(13, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("66÷5=13, 1", $true, $false, $false, $false, $false, $true, 0, $false, "65÷4=16, 1", 1) | Out-Null

$cell = $t.Cell(13, 5)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("56÷5=11, 1", $true, $false, $false, $false, $false, $true, 0, $false, "20÷9=2, 2", 1) | Out-Null

$cell = $t.Cell(17, 1)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("97÷7=13, 6", $true, $false, $false, $false, $false, $true, 0, $false, "64÷2=32, 0", 1) | Out-Null

$cell = $t.Cell(17, 2)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("56÷2=28, 0", $true, $false, $false, $false, $false, $true, 0, $false, "30÷7=4, 2", 1) | Out-Null

$cell = $t.Cell(17, 3)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("86÷6=14, 2", $true, $false, $false, $false, $false, $true, 0, $false, "48÷4=12, 0", 1) | Out-Null

$cell = $t.Cell(17, 4)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("57÷5=11, 2", $true, $false, $false, $false, $false, $true, 0, $false, "23÷7=3, 2", 1) | Out-Null

$cell = $t.Cell(17, 5)
$r = $d.Range($cell.Range.Start, $cell.Range.End)
$r.Find.Execute("56÷4=14, 0", $true, $false, $false, $false, $false, $true, 0, $false, "72÷7=10, 2", 1) | Out-Null
